# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme name "Office")    [only wired to the Notes Master]
#   ppt/theme/theme2.xml -> "Integral"     (clrScheme name "Red Violet")[wired to the Slide Master + the deck]
#
# The target edit swaps the *content* of those two parts: the theme actually
# driving the deck's look (the one reachable from SlideMaster/Slides) becomes
# the plain "Office" palette, while the Red Violet / Integral palette moves to
# the otherwise-unused part. fontScheme/fmtScheme are identical between the
# two themes already, so the only real difference is the 12 color-scheme
# entries (and the cosmetic name attributes, which aren't separately exposed
# on this object model).
#
# Apply it the supported way: write each of the 12 theme colors on the live
# ThemeColorScheme to the values the "Office" theme used, via .RGB (PowerPoint
# COM stores RGB as a single Long = R + G*256 + B*65536).

function Get-RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Index order matches the OOXML clrScheme child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    @(0x00, 0x00, 0x00),   # dk1
    @(0xFF, 0xFF, 0xFF),   # lt1
    @(0x44, 0x54, 0x6A),   # dk2
    @(0xE7, 0xE6, 0xE6),   # lt2
    @(0x5B, 0x9B, 0xD5),   # accent1
    @(0xED, 0x7D, 0x31),   # accent2
    @(0xA5, 0xA5, 0xA5),   # accent3
    @(0xFF, 0xC0, 0x00),   # accent4
    @(0x44, 0x72, 0xC4),   # accent5
    @(0x70, 0xAD, 0x47),   # accent6
    @(0x05, 0x63, 0xC1),   # hlink
    @(0x95, 0x4F, 0x72)    # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $c = $officeColors[$i]
    $tcs.Item($i + 1).RGB = Get-RGBVal $c[0] $c[1] $c[2]
}
